$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.978.99"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.81"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.57"
$ws.Range("E5").Value = "  +4.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5022"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2737"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06196"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07261"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.741.96"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6545"
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.18"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.728"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.80"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.999.73"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006859"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.617"
$ws.Range("E21").Value = "  +9.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.966.28"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.406"
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.58"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.494"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.28"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.790"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.51"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.000"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08167"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.705"
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04747"
$ws.Range("E33").Value = "  +4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.668"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9980"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6148"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.756"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01621"
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.936"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.13"
$ws.Range("E41").Value = "  +3.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8022"
$ws.Range("E42").Value = "  +9.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3921"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.013"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1176"
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.398"
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.91"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05295"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.682"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3489"
$ws.Range("E51").Value = "  +2.42%  "
